$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted as row 58 (existing rows 58-65
# shift down to 59-66, extending the used range to A1:R66).
$ws.Rows.Item(58).Insert()

$ws.Cells.Item(58, 1).Value = 11
$ws.Cells.Item(58, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(58, 3).Value = "Bíobío"
$ws.Cells.Item(58, 4).Value = 45077
$ws.Cells.Item(58, 5).Value = 8
$ws.Cells.Item(58, 6).Value = 100112030
$ws.Cells.Item(58, 7).Value = "Poroto granado"
$ws.Cells.Item(58, 8).Value = "Sin especificar"
$ws.Cells.Item(58, 9).Value = "Primera"
$ws.Cells.Item(58, 10).Value = 100
$ws.Cells.Item(58, 11).Value = 23000
$ws.Cells.Item(58, 12).Value = 24000
$ws.Cells.Item(58, 13).Value = 23500
$ws.Cells.Item(58, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(58, 15).Value = "Región Metropolitana"
$ws.Cells.Item(58, 16).Value = 940
$ws.Cells.Item(58, 17).Value = 25
$ws.Cells.Item(58, 18).Value = "Hortaliza"
